$wb = $excel.ActiveWorkbook

# ALC row 4: Root Rush
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 59.18182
$ws.Range("I4").Value = 46
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 46
$ws.Range("L4").Value = 75
$ws.Range("M4").Value = 68
$ws.Range("N4").Value = -303

# ARM row 31: I Was a Teenage Wailer
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 14738.667
$ws.Range("I31").Value = 5577.1665
$ws.Range("K31").Value = 5577.1665
$ws.Range("M31").Value = -5283.1665

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5434.905
$ws.Range("I32").Value = 4858.676
$ws.Range("J32").Value = 9699
$ws.Range("K32").Value = 4858.676
$ws.Range("L32").Value = 9699
$ws.Range("M32").Value = -4571.676
$ws.Range("N32").Value = -10273

# ARM row 45: Hollow Hallmarks
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3672.6365
$ws.Range("I45").Value = 3150.2
$ws.Range("J45").Value = 4108
$ws.Range("K45").Value = 3150.2
$ws.Range("L45").Value = 4108
$ws.Range("M45").Value = -2773.2
$ws.Range("N45").Value = -4862

# ARM row 63: Rivets Run through It
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2123.3076
$ws.Range("J63").Value = 3000
$ws.Range("L63").Value = 3000
$ws.Range("N63").Value = -4372

# ARM row 66: A Riveting Revival (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2123.3076
$ws.Range("J66").Value = 3000
$ws.Range("L66").Value = 15000
$ws.Range("N66").Value = -21864

# ARM row 102: Smells of Rich Tama-hagane
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 157109.61
$ws.Range("I102").Value = 184402.27
$ws.Range("J102").Value = 7000
$ws.Range("K102").Value = 184402.27
$ws.Range("L102").Value = 7000
$ws.Range("M102").Value = -182780.27
$ws.Range("N102").Value = -10244

# ARM row 128: Heading toward Bankruptcy
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 85000
$ws.Range("J128").Value = 85000
$ws.Range("L128").Value = 85000
$ws.Range("N128").Value = -94960

# BSM row 81: Diamond Sawdust
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0

# BSM row 82: Spirituality Inspector
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 29176.375
$ws.Range("I82").Value = 5744.75
$ws.Range("K82").Value = 5744.75
$ws.Range("M82").Value = -5361.75

# BSM row 84: I'm a Lumberjack and I'm Okay (L)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0

# BSM row 85: The Clamor for Hammers (L)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 29176.375
$ws.Range("I85").Value = 5744.75
$ws.Range("K85").Value = 5744.75
$ws.Range("M85").Value = -4418.75

# BSM row 132: Always Be Prepaired
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 93498
$ws.Range("J132").Value = 93498
$ws.Range("L132").Value = 93498
$ws.Range("N132").Value = -103618

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9263063
$ws.Range("J31").Value = 20838086
$ws.Range("L31").Value = 20838086
$ws.Range("N31").Value = -20838676

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9263063
$ws.Range("J34").Value = 20838086
$ws.Range("L34").Value = 20838086
$ws.Range("N34").Value = -20838490

# CRP row 58: You Do the Heavy Lifting
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1739.8
$ws.Range("I58").Value = 1879.8
$ws.Range("J58").Value = 1599.8
$ws.Range("K58").Value = 1879.8
$ws.Range("L58").Value = 1599.8
$ws.Range("M58").Value = -1676.8
$ws.Range("N58").Value = -2005.8

# CRP row 100: Run Before They Walk
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 51854.5
$ws.Range("I100").Value = 49709
$ws.Range("J100").Value = 54000
$ws.Range("K100").Value = 49709
$ws.Range("L100").Value = 54000
$ws.Range("M100").Value = -48627
$ws.Range("N100").Value = -56164

# CRP row 122: Timber of Tenkonto
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 26316112
$ws.Range("I122").Value = 649
$ws.Range("J122").Value = 52631576
$ws.Range("K122").Value = 1947
$ws.Range("L122").Value = 157894728
$ws.Range("M122").Value = 503
$ws.Range("N122").Value = -157899628

# CRP row 131: An Integral Reward
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 49466.332
$ws.Range("J131").Value = 49466.332
$ws.Range("L131").Value = 49466.332
$ws.Range("N131").Value = -59546.332

# CRP row 134: Wood You Be Quiet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2667.7942
$ws.Range("I134").Value = 2262.6667
$ws.Range("J134").Value = 4230.4287
$ws.Range("K134").Value = 6788.000100000001
$ws.Range("L134").Value = 12691.2861
$ws.Range("M134").Value = -4253.000100000001
$ws.Range("N134").Value = -17761.2861

# CRP row 136: Turali Quality
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1739.8
$ws.Range("I136").Value = 1879.8
$ws.Range("J136").Value = 1599.8
$ws.Range("K136").Value = 5639.4
$ws.Range("L136").Value = 4799.4
$ws.Range("M136").Value = -3089.4
$ws.Range("N136").Value = -9899.4

# CUL row 2: Pork Is a Salty Food
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 182.72728
$ws.Range("J2").Value = 20
$ws.Range("L2").Value = 120
$ws.Range("N2").Value = -346

# CUL row 92: Oh No Udon
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1549.8
$ws.Range("I92").Value = 1250
$ws.Range("J92").Value = 1999.5
$ws.Range("K92").Value = 3750
$ws.Range("L92").Value = 5998.5
$ws.Range("M92").Value = -2502
$ws.Range("N92").Value = -8494.5

# CUL row 111: Soup for the Soldier
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 559.5
$ws.Range("I111").Value = 559.5
$ws.Range("K111").Value = 1678.5
$ws.Range("M111").Value = 1388.5

# CUL row 113: Can't Eat Just One
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1197.75
$ws.Range("J113").Value = 2158
$ws.Range("L113").Value = 6474
$ws.Range("N113").Value = -10814

# CUL row 118: Teetotally
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 740.4
$ws.Range("I118").Value = 740.4
$ws.Range("K118").Value = 2221.2
$ws.Range("M118").Value = -978.1999999999998

# CUL row 121: A Cookie for Your Troubles
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1296.2858
$ws.Range("J121").Value = 893.75
$ws.Range("L121").Value = 2681.25
$ws.Range("N121").Value = -5301.25

# CUL row 122: Salt of the North
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2061.7827
$ws.Range("J122").Value = 2461.889
$ws.Range("L122").Value = 22157.001
$ws.Range("N122").Value = -27057.001

# CUL row 124: Bobbing for Compliments
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 47622028
$ws.Range("I124").Value = 3476.5
$ws.Range("K124").Value = 10429.5
$ws.Range("M124").Value = -5519.5

# CUL row 125: At Any Temperature
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 166669170
$ws.Range("J125").Value = 166669170
$ws.Range("L125").Value = 500007510
$ws.Range("N125").Value = -500017350

# CUL row 126: Imperial Palate
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 104172720
$ws.Range("I126").Value = 100004180
$ws.Range("J126").Value = 111120280
$ws.Range("K126").Value = 300012540
$ws.Range("L126").Value = 333360840
$ws.Range("M126").Value = -300007600
$ws.Range("N126").Value = -333370720

# CUL row 138: Bring Me Your Tacos
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 6880
$ws.Range("J138").Value = 3250
$ws.Range("L138").Value = 9750
$ws.Range("N138").Value = -20030

# CUL row 139: Najoothie
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2714.8
$ws.Range("I139").Value = 2635.25
$ws.Range("K139").Value = 7905.75
$ws.Range("M139").Value = -2765.75

# CUL row 140: Sweet, Sweet Bean Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2157.6667
$ws.Range("I140").Value = 1758.4445
$ws.Range("K140").Value = 5275.333500000001
$ws.Range("M140").Value = -95.33350000000064

# GSM row 55: If You've Got It, Flaunt It
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 10000
$ws.Range("K55").Value = 10000
$ws.Range("M55").Value = -9673

# GSM row 126: Gold Rush Order
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2046.1428
$ws.Range("I126").Value = 1887.1666
$ws.Range("K126").Value = 5661.4998
$ws.Range("M126").Value = -3191.4998

# GSM row 134: Guaranteed Gem
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 57855.2
$ws.Range("J134").Value = 57855.2
$ws.Range("L134").Value = 173565.6
$ws.Range("N134").Value = -178635.6

# LTW row 16: Saddle Sore
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 952.1667
$ws.Range("J16").Value = 2002
$ws.Range("L16").Value = 2002
$ws.Range("N16").Value = -2342

# LTW row 102: Shrug It On
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 63500
$ws.Range("I102").Value = 63500
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 63500
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -60255

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4076.4614
$ws.Range("I136").Value = 2799.4
$ws.Range("K136").Value = 8398.200000000001
$ws.Range("M136").Value = -5848.200000000001

# WVR row 102: Don't Sweat the Role
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 69000
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# WVR row 123: Helping Handwear
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 74103.22
$ws.Range("J123").Value = 72053.625
$ws.Range("L123").Value = 72053.625
$ws.Range("N123").Value = -81853.625
